$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.249.01'
$ws.Range('E2').Value = '  -1.91%  '
$ws.Range('D3').Value = '2.581.81'
$ws.Range('E3').Value = '  -2.96%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '563.32'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.09%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '142.59'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.40%  '
$ws.Range('E7').Value = '  +0.16%  '
$ws.Range('E8').Value = '  -2.13%  '
$ws.Range('D9').Value = '2.588.07'
$ws.Range('E9').Value = '  -2.65%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.62'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.85%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.103'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.47%  '
$ws.Range('E12').Value = '  +11.39%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.352'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.95%  '
$ws.Range('D14').Value = '3.037.36'
$ws.Range('E14').Value = '  -2.23%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '23.33'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +7.41%  '
$ws.Range('D16').Value = '59.202.37'
$ws.Range('E16').Value = '  -1.91%  '
$ws.Range('E17').Value = '  +0.41%  '
$ws.Range('D18').Value = '2.590.73'
$ws.Range('E18').Value = '  -2.38%  '
$ws.Range('E19').Value = '  +0.51%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '337.14'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.14%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.35'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.42%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.38'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.10%  '
$ws.Range('E23').Value = '  +0.26%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '64.14'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -4.03%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.465'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +5.41%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.998'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.45%  '
$ws.Range('E27').Value = '  -3.05%  '
$ws.Range('E28').Value = '  -0.31%  '
$ws.Range('E29').Value = '  -0.02%  '
$ws.Range('E30').Value = '  +0.00%  '
$ws.Range('E31').Value = '  -2.55%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '160.32'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.58%  '
$ws.Range('E33').Value = '  -0.37%  '
$ws.Range('E34').Value = '  -1.20%  '
$ws.Range('E35').Value = '  -1.47%  '
$ws.Range('E36').Value = '  -0.83%  '
$ws.Range('B37').Value = 'Fetch.AI'
$ws.Range('C37').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.879'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.34%  '
$ws.Range('B38').Value = 'SuiNetwork'
$ws.Range('C38').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.872'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.88%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '37.41'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.45%  '
$ws.Range('E40').Value = '  -2.25%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '294.41'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -4.72%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.67'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.03%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.999'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.21%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '131.79'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +4.74%  '
$ws.Range('E45').Value = '  -1.00%  '
$ws.Range('E46').Value = '  -1.81%  '
$ws.Range('B47').Value = 'WhiteBITCoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.65'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.06%  '
$ws.Range('B48').Value = 'Hedera'
$ws.Range('C48').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0535'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.25%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '18.99'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.93%  '
$ws.Range('E50').Value = '  -0.58%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '18.63'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.92%  '
